# Insert a new data row before the existing row 175 ("Ají" / Terminal La
# Palmera de La Serena log). This shifts the existing rows 175-204 down to
# 176-205 (carrying their formatting/styles along), matching the diff, which
# shows every row from 176 to 205 taking on the values previously held by
# the row above it, and a brand-new record landing in row 175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new record's data.
$ws.Cells.Item(175, 1).Value = 8
$ws.Cells.Item(175, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(175, 3).Value = "Coquimbo"
$ws.Cells.Item(175, 4).Value = 44641
$ws.Cells.Item(175, 5).Value = 4
$ws.Cells.Item(175, 6).Value = 100112021
$ws.Cells.Item(175, 7).Value = "Ají"
$ws.Cells.Item(175, 8).Value = "Inferno"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 400
$ws.Cells.Item(175, 11).Value = 17000
$ws.Cells.Item(175, 12).Value = 18000
$ws.Cells.Item(175, 13).Value = 17500
$ws.Cells.Item(175, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(175, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(175, 16).Value = 1167
$ws.Cells.Item(175, 17).Value = 15
$ws.Cells.Item(175, 18).Value = "Hortaliza"
